# "some old code move to other folder for now" — drop the "Mars" and "A13"/"A11"
# sample rows from the VfM tables (sheets Q1_20_21 and Q4_19_20), re-order the
# remaining Q4_19_20 rows (Columbia before F9), and refresh the Count sheet's
# roll-up totals/counts that fed off the removed rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Q1_20_21": delete the "Mars" row and the "A13" row.
# Before: Mars(3), SoT(4), A13(5), F9(6), Columbia(7)
# After:  SoT(3), F9(4), Columbia(5)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Q1_20_21")
$ws1.Rows.Item(3).Delete()   # removes "Mars" (SoT shifts up to row 3)
$ws1.Rows.Item(4).Delete()   # removes "A13" (now at row 4 after first delete)

# ---------------------------------------------------------------------------
# Sheet "Q4_19_20": delete the "Mars" row and the near-empty "A11" row, then
# swap the remaining F9/Columbia rows so Columbia comes before F9.
# Before: Mars(3), SoT(4), A11(5), A13(6), F9(7), Columbia(8)
# After delete: SoT(3), A13(4), F9(5), Columbia(6)
# After swap:   SoT(3), A13(4), Columbia(5), F9(6)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Q4_19_20")
$ws2.Rows.Item(3).Delete()   # removes "Mars" (SoT shifts up to row 3)
$ws2.Rows.Item(4).Delete()   # removes "A11" (A13 shifts up to row 4)

$rowF9 = $ws2.Range("B5:K5")
$rowColumbia = $ws2.Range("B6:K6")
$f9Values = $rowF9.Value2
$columbiaValues = $rowColumbia.Value2
$rowF9.Value2 = $columbiaValues
$rowColumbia.Value2 = $f9Values

# ---------------------------------------------------------------------------
# Sheet "Count": refresh the PVC totals and category counts so they match the
# trimmed data above (Mars's "Very High" PVC and A13's "High" PVC are gone,
# and the Q1/Q4 totals + counts drop accordingly).
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Count")

$ws3.Range("C7").Value2 = 0      # High / Q1 20-21 PVC total
$ws3.Range("C8").Value2 = 0      # Very High / Q1 20-21 PVC total
$ws3.Range("D8").Value2 = 0      # Very High / Q4 19-20 PVC total
$ws3.Range("C11").Value2 = 4128  # Total / Q1 20-21 PVC total
$ws3.Range("D11").Value2 = 5764  # Total / Q4 19-20 PVC total

$ws3.Range("C19").Value2 = 1     # High / Q1 20-21 count
$ws3.Range("C20").Value2 = 0     # Very High / Q1 20-21 count
$ws3.Range("D20").Value2 = 0     # Very High / Q4 19-20 count
$ws3.Range("C23").Value2 = 3     # Total / Q1 20-21 count
$ws3.Range("D23").Value2 = 4     # Total / Q4 19-20 count
